$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I2").Value = 0.2594103048008066
$ws.Range("J2").Value = 0.2594103048008066
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 1.949330262895334
$ws.Range("R2").Value = 17.543972366058
$ws.Range("S2").Value = 0.07434537852950945
$ws.Range("T2").Value = 0.07434537852950944
$ws.Range("I3").Value = 0.2594103048008066
$ws.Range("J3").Value = 0.2594103048008066
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("S3").Value = 0.08613187131701081
$ws.Range("T3").Value = 0.08613187131701078
$ws.Range("I4").Value = 0.2594103048008066
$ws.Range("J4").Value = 0.2594103048008066
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 1.714199823285333
$ws.Range("R4").Value = 15.427798409568
$ws.Range("S4").Value = 0.06537775417700434
$ws.Range("T4").Value = 0.06537775417700432
$ws.Range("I5").Value = 0.2594103048008066
$ws.Range("J5").Value = 0.2594103048008066
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 0.8798174759409999
$ws.Range("R5").Value = 7.918357283468999
$ws.Range("S5").Value = 0.03355530077728205
$ws.Range("T5").Value = 0.03355530077728205
$ws.Range("G6").Value = 0.03819566666666666
$ws.Range("H6").Value = 0.114587
$ws.Range("I6").Value = 0.7405896951991934
$ws.Range("J6").Value = 0.7405896951991934
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 5.565137076373111
$ws.Range("R6").Value = 50.086233687358
$ws.Range("S6").Value = 0.2122483964810748
$ws.Range("T6").Value = 0.2122483964810748
$ws.Range("G7").Value = 0.03819566666666666
$ws.Range("H7").Value = 0.114587
$ws.Range("I7").Value = 0.7405896951991934
$ws.Range("J7").Value = 0.7405896951991934
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("Q7").Value = 6.447417176488445
$ws.Range("R7").Value = 58.026754588396
$ws.Range("S7").Value = 0.2458976191195733
$ws.Range("T7").Value = 0.2458976191195733
$ws.Range("G8").Value = 0.03819566666666666
$ws.Range("H8").Value = 0.114587
$ws.Range("I8").Value = 0.7405896951991934
$ws.Range("J8").Value = 0.7405896951991934
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 4.893863894929778
$ws.Range("R8").Value = 44.044775054368
$ws.Range("S8").Value = 0.1866467528186061
$ws.Range("T8").Value = 0.1866467528186061
$ws.Range("G9").Value = 0.03819566666666666
$ws.Range("H9").Value = 0.114587
$ws.Range("I9").Value = 0.7405896951991934
$ws.Range("J9").Value = 0.7405896951991934
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 2.511788253124333
$ws.Range("R9").Value = 22.606094278119
$ws.Range("S9").Value = 0.09579692677993917
$ws.Range("T9").Value = 0.09579692677993917
